$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Compras": fill rows 6-11 (items 11-16) with purchase data and link
# each item number back to its corresponding row in Inventario_cero_coste
# (row 8's item links to a row within the Compras sheet itself).
# ---------------------------------------------------------------------------
$wsCompras = $wb.Worksheets.Item("Compras")

$comprasRows = @(
    @{ Row = 6;  Num = 11; Precio = 22.49; Target = "Inventario_cero_coste!B13"; Display = "11" },
    @{ Row = 7;  Num = 12; Precio = 13.99; Target = "Inventario_cero_coste!B14"; Display = "12" },
    @{ Row = 8;  Num = 13; Precio = 6.95;  Target = "Compras!B15";               Display = "13" },
    @{ Row = 9;  Num = 14; Precio = 6.95;  Target = "Inventario_cero_coste!B16"; Display = "15" },
    @{ Row = 10; Num = 15; Precio = 8.99;  Target = "Inventario_cero_coste!B17"; Display = "15" },
    @{ Row = 11; Num = 16; Precio = 22.09; Target = "Inventario_cero_coste!B18"; Display = "16" }
)

foreach ($item in $comprasRows) {
    $r = $item.Row
    # Add the hyperlink first (it sets a text caption), then overwrite the
    # cell with its real numeric item value so the stored value stays a
    # number while the hyperlink keeps its own "display" caption text.
    $wsCompras.Hyperlinks.Add($wsCompras.Range("B$r"), "", $item.Target, [Type]::Missing, $item.Display) | Out-Null
    $wsCompras.Range("B$r").Value = $item.Num
    $wsCompras.Range("C$r").Value = "Manu"
    $wsCompras.Range("D$r").Value = $item.Precio
}

# ---------------------------------------------------------------------------
# Sheet "Inventario_cero_coste": fill rows 13-18 (items 11-16) with the new
# mechanical / electronic components, their descriptions and purchase links.
# (Row 13's "F" note is filled in last, mirroring the original authoring
# order so the shared-string table layout matches.)
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("Inventario_cero_coste")

$wsInv.Range("C13").Value = "DRV8825 "
$wsInv.Range("E13").Value = "Módulo controlador Stepper"
$wsInv.Range("G13").Value = "Manu"
$wsInv.Range("I13").Value = "https://www.amazon.es/dp/B07YWV6W4W?psc=1&ref=ppx_yo2ov_dt_b_product_details"
$wsInv.Hyperlinks.Add($wsInv.Range("I13"), "https://www.amazon.es/dp/B07YWV6W4W?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", [Type]::Missing, "https://www.amazon.es/dp/B07YWV6W4W?psc=1&ref=ppx_yo2ov_dt_b_product_details") | Out-Null

$wsInv.Range("C14").Value = "Bolas rodamientos"
$wsInv.Range("E14").Value = "Bolas rodamientos"
$wsInv.Range("F14").Value = "Bolas rodamientos 560 piezas, 11 tamaños"
$wsInv.Range("G14").Value = "Manu"
$wsInv.Range("I14").Value = "https://www.amazon.es/dp/B094346M9W?psc=1&ref=ppx_yo2ov_dt_b_product_details"
$wsInv.Hyperlinks.Add($wsInv.Range("I14"), "https://www.amazon.es/dp/B094346M9W?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", [Type]::Missing, "https://www.amazon.es/dp/B094346M9W?psc=1&ref=ppx_yo2ov_dt_b_product_details") | Out-Null

$wsInv.Range("C15").Value = "Hilo tórico"
$wsInv.Range("E15").Value = "Hilo tórico"
$wsInv.Range("F15").Value = "Hilo tórico de caucho 2.5 mm de diámetro, 1 m"
$wsInv.Range("G15").Value = "Manu"
$wsInv.Range("I15").Value = "https://www.amazon.es/dp/B0BD5P49Q8?psc=1&ref=ppx_yo2ov_dt_b_product_details"
$wsInv.Hyperlinks.Add($wsInv.Range("I15"), "https://www.amazon.es/dp/B0BD5P49Q8?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", [Type]::Missing, "https://www.amazon.es/dp/B0BD5P49Q8?psc=1&ref=ppx_yo2ov_dt_b_product_details") | Out-Null

$wsInv.Range("C16").Value = "Hilo tórico"
$wsInv.Range("E16").Value = "Hilo tórico"
$wsInv.Range("F16").Value = "Hilo tórico de caucho 2.5 mm de diámetro, 1 m"
$wsInv.Range("G16").Value = "Manu"
$wsInv.Range("I16").Value = "https://www.amazon.es/dp/B0BD5P49Q8?psc=1&ref=ppx_yo2ov_dt_b_product_details"
$wsInv.Hyperlinks.Add($wsInv.Range("I16"), "https://www.amazon.es/dp/B0BD5P49Q8?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", [Type]::Missing, "https://www.amazon.es/dp/B0BD5P49Q8?psc=1&ref=ppx_yo2ov_dt_b_product_details") | Out-Null

$wsInv.Range("G17").Value = "Manu"
$wsInv.Range("I17").Value = "https://www.amazon.es/dp/B0BB15P6DL?psc=1&ref=ppx_yo2ov_dt_b_product_details"
$wsInv.Hyperlinks.Add($wsInv.Range("I17"), "https://www.amazon.es/dp/B0BB15P6DL?psc=1&ref=ppx_yo2ov_dt_b_product_details", "", [Type]::Missing, "https://www.amazon.es/dp/B0BB15P6DL?psc=1&ref=ppx_yo2ov_dt_b_product_details") | Out-Null
$wsInv.Range("C17").Value = "Cubrecables"
$wsInv.Range("E17").Value = "Cubrecables"
$wsInv.Range("F17").Value = "Cubrecables 13-20 mm de diámetro, 3 m"

$wsInv.Range("C18").Value = "Mecánica"
$wsInv.Range("E18").Value = "Elementos mecánicos"
$wsInv.Range("F18").Value = "2 Correas dentadas GT2, 5 rodamientos lineales, 2 varillas lisas, 3 poleas GT2"
$wsInv.Range("G18").Value = "Manu"
$wsInv.Range("I18").Value = "https://www.hta3d.com"
$wsInv.Hyperlinks.Add($wsInv.Range("I18"), "https://www.hta3d.com", "", [Type]::Missing, "https://www.hta3d.com") | Out-Null

$wsInv.Range("F13").Value = "5 Módulos controlador para stepper"

# ---------------------------------------------------------------------------
# Update sheet selections / active sheet: move focus from Inventario back to
# Compras, keep the last selection used on Inventario (row 13 area).
# ---------------------------------------------------------------------------
$wsInv.Range("B13").Select() | Out-Null
$wsCompras.Activate() | Out-Null
$wsCompras.Range("B11").Select() | Out-Null
